$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.056054599942534658
$ws.Range("B1").Value = -0.056263324147086652
$ws.Range("A2").Value = -0.060019046902855672
$ws.Range("B2").Value = -0.060019425050125252
$ws.Range("A3").Value = -0.024621027669214476
$ws.Range("B3").Value = -0.024615383348255516
